$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The XPath locator strings used in row 1 (C1:G1) switched from double
# quotes to single quotes around the attribute value. Re-writing the
# cell values causes the old (now-unused) shared-string entries to be
# dropped and the new text appended to the shared-string table, which
# is exactly what the target file does.
$ws.Range("C1").Value = "//*[@id='insurance-form']/div/section[1]"
$ws.Range("D1").Value = "//*[@id='insurance-form']/div/section[2]"
$ws.Range("E1").Value = "//*[@id='insurance-form']/div/section[3]"
$ws.Range("F1").Value = "//*[@id='insurance-form']/div/section[4]"
$ws.Range("G1").Value = "//*[@id='insurance-form']/div/section[5]"

# Move the sheet selection from A19:XFD24 to G11.
$ws.Range("G11").Select()
